$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values scraped on Tue Sep 26 10:36:37 UTC 2023
$updates = @(
    @{Cell='D2'; Value='26.259.41'}
    @{Cell='E2'; Value='  +0.44%  '}
    @{Cell='D3'; Value='1.588.38'}
    @{Cell='E3'; Value='  +0.67%  '}
    @{Cell='E4'; Value='  -0.16%  '}
    @{Cell='D5'; Value='212.39'}
    @{Cell='E5'; Value='  +1.58%  '}
    @{Cell='D6'; Value='0.500'}
    @{Cell='E6'; Value='  +0.54%  '}
    @{Cell='E7'; Value='  -0.16%  '}
    @{Cell='E8'; Value='  +0.26%  '}
    @{Cell='E9'; Value='  -0.19%  '}
    @{Cell='D10'; Value='19.33'}
    @{Cell='E10'; Value='  -0.90%  '}
    @{Cell='D11'; Value='0.0850'}
    @{Cell='E11'; Value='  +0.86%  '}
    @{Cell='D12'; Value='1.813.04'}
    @{Cell='E12'; Value='  +0.72%  '}
    @{Cell='D13'; Value='1.595.14'}
    @{Cell='E13'; Value='  +1.38%  '}
    @{Cell='E14'; Value='  -0.10%  '}
    @{Cell='E15'; Value='  +1.27%  '}
    @{Cell='D16'; Value='64.32'}
    @{Cell='E16'; Value='  -0.13%  '}
    @{Cell='D17'; Value='26.247.09'}
    @{Cell='E17'; Value='  +0.37%  '}
    @{Cell='E18'; Value='  -0.40%  '}
    @{Cell='E19'; Value='  +2.30%  '}
    @{Cell='D20'; Value='212.54'}
    @{Cell='E20'; Value='  +2.33%  '}
    @{Cell='E21'; Value='  -0.11%  '}
    @{Cell='E22'; Value='  +1.09%  '}
    @{Cell='E23'; Value='  +1.52%  '}
    @{Cell='E24'; Value='  -2.73%  '}
    @{Cell='D25'; Value='143.86'}
    @{Cell='E25'; Value='  +0.09%  '}
    @{Cell='E26'; Value='  -0.15%  '}
    @{Cell='E27'; Value='  +1.38%  '}
    @{Cell='E28'; Value='  -0.37%  '}
    @{Cell='D29'; Value='15.19'}
    @{Cell='E29'; Value='  -0.13%  '}
    @{Cell='E30'; Value='  -1.46%  '}
    @{Cell='E31'; Value='  +1.04%  '}
    @{Cell='E32'; Value='  -0.20%  '}
    @{Cell='E33'; Value='  -0.96%  '}
    @{Cell='D34'; Value='1.334.24'}
    @{Cell='E34'; Value='  +4.36%  '}
    @{Cell='E35'; Value='  -0.93%  '}
    @{Cell='E36'; Value='  -0.64%  '}
    @{Cell='E37'; Value='  -3.42%  '}
    @{Cell='E38'; Value='  +0.41%  '}
    @{Cell='D39'; Value='0.817'}
    @{Cell='E39'; Value='  +0.25%  '}
    @{Cell='E40'; Value='  -6.44%  '}
    @{Cell='D41'; Value='5.72'}
    @{Cell='E41'; Value='  +2.94%  '}
    @{Cell='E42'; Value='  -0.12%  '}
    @{Cell='E43'; Value='  +0.08%  '}
    @{Cell='D44'; Value='0.762'}
    @{Cell='E44'; Value='  -0.04%  '}
    @{Cell='B45'; Value='Aave'}
    @{Cell='C45'; Value='https://coinranking.com/coin/ixgUfzmLR+aave-aave'}
    @{Cell='D45'; Value='61.84'}
    @{Cell='E45'; Value='  -0.78%  '}
    @{Cell='B46'; Value='RocketPoolETH'}
    @{Cell='C46'; Value='https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'}
    @{Cell='D46'; Value='1.724.58'}
    @{Cell='E46'; Value='  +0.65%  '}
    @{Cell='D47'; Value='85.64'}
    @{Cell='E47'; Value='  -3.59%  '}
    @{Cell='E48'; Value='  -3.80%  '}
    @{Cell='E49'; Value='  -0.78%  '}
    @{Cell='E50'; Value='  -2.88%  '}
    @{Cell='E51'; Value='  -0.25%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
